# Apply the crypto price/volume refresh described in the commit diff.
# Numeric-looking Price (column D) values are prefixed with a literal
# leading apostrophe so Excel stores them as text (e.g. "1.000"), matching
# the original inline-string cell contents instead of coercing to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.730.99"
$ws.Range("E2").Value = "  -1.85%  "

$ws.Range("D3").Value = "1.867.84"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'300.60"
$ws.Range("E5").Value = "  -2.19%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'0.5317"
$ws.Range("E7").Value = "  +1.08%  "

$ws.Range("D8").Value = "'0.3725"
$ws.Range("E8").Value = "  -2.33%  "

$ws.Range("E9").Value = "  -1.54%  "

$ws.Range("D10").Value = "'21.44"
$ws.Range("E10").Value = "  -1.55%  "

$ws.Range("D11").Value = "'0.8872"
$ws.Range("E11").Value = "  -1.74%  "

$ws.Range("D12").Value = "'0.08181"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "1.875.60"
$ws.Range("E13").Value = "  +28.93%  "

$ws.Range("D14").Value = "'92.29"
$ws.Range("E14").Value = "  -4.13%  "

$ws.Range("D15").Value = "'5.289"
$ws.Range("E15").Value = "  -1.37%  "

$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").Value = "'14.77"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("D18").Value = "'0.000008481"
$ws.Range("E18").Value = "  -2.07%  "

$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "26.760.53"
$ws.Range("E20").Value = "  -1.89%  "

$ws.Range("D21").Value = "'4.971"
$ws.Range("E21").Value = "  -2.89%  "

$ws.Range("D22").Value = "'10.60"
$ws.Range("E22").Value = "  -2.11%  "

$ws.Range("D23").Value = "'6.355"
$ws.Range("E23").Value = "  -2.40%  "

$ws.Range("D24").Value = "'2.285"
$ws.Range("E24").Value = "  -1.05%  "

$ws.Range("D25").Value = "'145.58"
$ws.Range("E25").Value = "  -3.05%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'18.00"
$ws.Range("E26").Value = "  -1.38%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'1.722"
$ws.Range("E27").Value = "  -1.21%  "

$ws.Range("D28").Value = "'113.39"
$ws.Range("E28").Value = "  -2.85%  "

$ws.Range("E29").Value = "  -3.35%  "

$ws.Range("E30").Value = "  -4.78%  "

$ws.Range("D31").Value = "'0.09118"
$ws.Range("E31").Value = "  -1.36%  "

$ws.Range("D32").Value = "'0.7984"
$ws.Range("E32").Value = "  -4.37%  "

$ws.Range("D33").Value = "'0.05004"
$ws.Range("E33").Value = "  -1.14%  "

$ws.Range("D34").Value = "'1.171"
$ws.Range("E34").Value = "  -4.65%  "

$ws.Range("D35").Value = "'2.944"
$ws.Range("E35").Value = "  -1.56%  "

$ws.Range("D36").Value = "'0.6055"
$ws.Range("E36").Value = "  +4.01%  "

$ws.Range("D37").Value = "'2.670"
$ws.Range("E37").Value = "  -2.06%  "

$ws.Range("D38").Value = "'3.170"
$ws.Range("E38").Value = "  -5.09%  "

$ws.Range("D39").Value = "'0.01942"
$ws.Range("E39").Value = "  -3.20%  "

$ws.Range("E40").Value = "  -1.41%  "

$ws.Range("D41").Value = "'6.491"
$ws.Range("E41").Value = "  -1.71%  "

$ws.Range("D42").Value = "'0.5192"
$ws.Range("E42").Value = "  +5.26%  "

$ws.Range("D43").Value = "'8.707"
$ws.Range("E43").Value = "  -5.07%  "

$ws.Range("D44").Value = "'114.60"
$ws.Range("E44").Value = "  -2.45%  "

$ws.Range("D45").Value = "'0.1489"
$ws.Range("E45").Value = "  -2.06%  "

$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").Value = "'9.952"
$ws.Range("E47").Value = "  -2.02%  "

$ws.Range("D48").Value = "'1.632"
$ws.Range("E48").Value = "  -0.65%  "

$ws.Range("D49").Value = "'37.40"
$ws.Range("E49").Value = "  -3.90%  "

$ws.Range("D50").Value = "'0.06043"
$ws.Range("E50").Value = "  -1.56%  "

$ws.Range("D51").Value = "'61.98"
$ws.Range("E51").Value = "  -3.86%  "
